$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.014.96"
$ws.Range("E2").Value = "  -2.41%  "
$ws.Range("D3").Value = "2.219.69"
$ws.Range("E3").Value = "  -4.32%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'244.45"
$ws.Range("E5").Value = "  -3.64%  "
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("D7").Value = "'73.81"
$ws.Range("E7").Value = "  -3.24%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -6.28%  "
$ws.Range("D10").Value = "'40.35"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D11").Value = "'0.0935"
$ws.Range("E11").Value = "  -5.39%  "
$ws.Range("D12").Value = "'7.01"
$ws.Range("E12").Value = "  -7.74%  "
$ws.Range("E13").Value = "  -3.59%  "
$ws.Range("D14").Value = "2.547.15"
$ws.Range("E14").Value = "  -4.67%  "
$ws.Range("D15").Value = "'14.40"
$ws.Range("E15").Value = "  -6.88%  "
$ws.Range("D16").Value = "'0.844"
$ws.Range("E16").Value = "  -4.93%  "
$ws.Range("D17").Value = "2.223.31"
$ws.Range("E17").Value = "  -4.64%  "
$ws.Range("D18").Value = "41.812.25"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("D19").Value = "0.0₃0964"
$ws.Range("E19").Value = "  -4.59%  "
$ws.Range("D20").Value = "'70.93"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("E21").Value = "  -5.36%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "'228.15"
$ws.Range("E23").Value = "  -4.42%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -5.99%  "
$ws.Range("D26").Value = "'10.96"
$ws.Range("E26").Value = "  -5.87%  "
$ws.Range("D27").Value = "'2.25"
$ws.Range("E27").Value = "  -7.50%  "
$ws.Range("D28").Value = "'7.18"
$ws.Range("E28").Value = "  +14.17%  "
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("D30").Value = "'167.68"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "'20.33"
$ws.Range("E31").Value = "  -4.71%  "
$ws.Range("D32").Value = "'0.0808"
$ws.Range("E32").Value = "  -4.75%  "
$ws.Range("D33").Value = "'30.34"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  -9.10%  "
$ws.Range("E35").Value = "  -3.33%  "
$ws.Range("E36").Value = "  -5.45%  "
$ws.Range("D37").Value = "'4.78"
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("D38").Value = "'0.0294"
$ws.Range("E38").Value = "  -6.77%  "
$ws.Range("D39").Value = "'13.15"
$ws.Range("E39").Value = "  -7.26%  "
$ws.Range("E40").Value = "  -9.19%  "
$ws.Range("E41").Value = "  -3.88%  "
$ws.Range("D42").Value = "'110.25"
$ws.Range("E42").Value = "  +3.56%  "
$ws.Range("E43").Value = "  -9.28%  "
$ws.Range("D44").Value = "'59.38"
$ws.Range("E44").Value = "  -5.51%  "
$ws.Range("D45").Value = "'8.58"
$ws.Range("E45").Value = "  -7.02%  "
$ws.Range("D46").Value = "'0.0984"
$ws.Range("E46").Value = "  -4.40%  "
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("E48").Value = "  -6.02%  "
$ws.Range("E49").Value = "  -3.61%  "
$ws.Range("D50").Value = "'4.15"
$ws.Range("E50").Value = "  -15.50%  "
$ws.Range("E51").Value = "  -1.31%  "
